$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values
$ws.Range("B2").Value = 262.52707065191231
$ws.Range("C2").Value = 303.32534830605175
$ws.Range("D2").Value = 257.34769713850034
$ws.Range("E2").Value = 308.60292676760866

# Row 3 data values
$ws.Range("B3").Value = 250.20360560472221
$ws.Range("C3").Value = 303.99727452461184
$ws.Range("D3").Value = 251.96893604029032
$ws.Range("E3").Value = 308.59097660620898

# Update selection to match new diff (B1:E3)
$ws.Range("B1:E3").Select()
